$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at row 132 (pushes the former row 132..169
# down to 133..170, matching the diff's observed shift of all subsequent
# rows).
$ws.Rows.Item(132).Insert()

$ws.Cells.Item(132, 1).Value = 10
$ws.Cells.Item(132, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(132, 3).Value = "La Araucanía"
$ws.Cells.Item(132, 4).Value = 44876
$ws.Cells.Item(132, 5).Value = 9
$ws.Cells.Item(132, 6).Value = "Fruta"
$ws.Cells.Item(132, 7).Value = 100107
$ws.Cells.Item(132, 8).Value = "Otros"
$ws.Cells.Item(132, 9).Value = 100107002
$ws.Cells.Item(132, 10).Value = "Chirimoya"
$ws.Cells.Item(132, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(132, 12).Value = "Primera"
$ws.Cells.Item(132, 13).Value = 40
$ws.Cells.Item(132, 14).Value = 2800
$ws.Cells.Item(132, 15).Value = 3000
$ws.Cells.Item(132, 16).Value = 2900
$ws.Cells.Item(132, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(132, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(132, 19).Value = 2900
$ws.Cells.Item(132, 20).Value = 1
